# Refresh the "Open-OMs" list with the latest OM numbers pulled from the
# PDF readers (see commit message). The previous 15 closed OMs are replaced
# by 32 new OMs; the Status column is cleared since these are not yet
# "Encerrado!" (closed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newOMs = @(
    "5197072",
    "5197173",
    "5196791",
    "5197351",
    "5197174",
    "5196603",
    "5196895",
    "5196780",
    "5197447",
    "5196812",
    "5197114",
    "5197572",
    "685601370003",
    "685601380872",
    "685601374545",
    "685601378879",
    "685601378435",
    "685601378437",
    "685601377298",
    "685601382621",
    "685601381340",
    "685601367925",
    "685601378723",
    "685601378307",
    "685601383769",
    "685601355769",
    "685601375759",
    "685601375604",
    "685601370724",
    "685601378491",
    "685601364160",
    "685601377899"
)

$firstDataRow = 2
$lastDataRow = $firstDataRow + $newOMs.Length - 1

# The OM numbers are stored as text (they come straight out of the PDF
# text-extraction step), so force the column to Text before writing the
# values -- otherwise the long digit strings would be auto-coerced to
# numbers.
$ws.Range("A$($firstDataRow):A$($lastDataRow)").NumberFormat = "@"

for ($i = 0; $i -lt $newOMs.Length; $i++) {
    $row = $firstDataRow + $i
    $ws.Cells.Item($row, 1).Value = $newOMs[$i]
    $ws.Cells.Item($row, 2).ClearContents()
}

Write-Output "Wrote $($newOMs.Length) OM rows ($firstDataRow to $lastDataRow)"
